# Update countries & provincias Spain
# - Refresh the "Datos actualizados" timestamp in A1
# - Apply the latest COVID-19 country counts, which re-sorts a few
#   neighbouring rows (the country that received fresh numbers overtakes
#   the country below it, which keeps its previous totals).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CountryRow($Row, $Country, $Total, $Nuevos, $Activos, $Recuperados, $Criticos, $MuertesHoy, $Muertes) {
    $ws.Cells.Item($Row, 1).Value = $Country
    $ws.Cells.Item($Row, 2).Value = $Total
    $ws.Cells.Item($Row, 3).Value = $Nuevos
    $ws.Cells.Item($Row, 4).Value = $Activos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $Criticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# --- Header / timestamp -----------------------------------------------
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 21 de Abril de 2020 a las 08:52"

# --- Row 33: Rumania (no reorder, counts updated) ------------------------
Set-CountryRow 33 "Rumania" 8936 0 2017 6437 261 4 482

# --- Row 39: Chequia (no reorder, counts updated) -------------------------
Set-CountryRow 39 "Chequia" 6914 14 1597 5121 75 2 196

# --- Rows 45-46: Ucrania overtakes Catar ----------------------------------
Set-CountryRow 45 "Ucrania" 6125 415 367 5597 45 10 161
Set-CountryRow 46 "Catar" 6015 0 555 5451 37 0 9

# --- Row 64: Kazajistan (no reorder, counts updated) ----------------------
Set-CountryRow 64 "Kazajistan" 1967 115 452 1496 22 0 19

# --- Row 103: Honduras (no reorder, counts updated) -----------------------
Set-CountryRow 103 "Honduras" 494 17 29 419 10 0 46

# --- Rows 107-108: Taiwan overtakes Jordania ------------------------------
Set-CountryRow 107 "Taiwan" 425 3 217 202 0 0 6
Set-CountryRow 108 "Jordania" 425 0 282 136 5 0 7

# --- Rows 125-126: El Salvador overtakes Jamaica --------------------------
Set-CountryRow 125 "El Salvador" 225 7 48 170 2 0 7
Set-CountryRow 126 "Jamaica" 223 0 27 191 0 0 5

# --- Row 128: Islas Feroe (no reorder, counts updated) --------------------
Set-CountryRow 128 "Islas Feroe" 185 0 178 7 0 0 0
